# Adapt model interface for addition of lever starts and ends.
#
# The "Control" sheet currently drives the model via a single
# Ambition column (B): INPUT.FACTOR = Control!$B$1, INPUT.OFFSET = Control!$B$2,
# referenced by named range input.lever.ambition (Control!$B$1:$B$2).
#
# This adds a "start" (C) and "end" (D) column to Control, exposed as two
# new named ranges (input.lever.start / input.lever.end), and wires the
# Sheet1 output formulas to take the (start - end) delta into account.

$wb = $excel.ActiveWorkbook

$control = $wb.Worksheets.Item("Control")

# New "start" / "end" lever columns alongside the existing Ambition column.
$control.Range("C1").Value = 1
$control.Range("D1").Value = 1
$control.Range("C2").Value = 2
$control.Range("D2").Value = 2

# Expose them as named ranges, matching the existing input.lever.ambition style.
$wb.Names.Add("input.lever.start", "=Control!`$C`$1:`$C`$2")
$wb.Names.Add("input.lever.end", "=Control!`$D`$1:`$D`$2")

# Wire the new start/end lever columns into the Sheet1 output rows.
$sheet1 = $wb.Worksheets.Item("Sheet1")

for ($col = 1; $col -le 21; $col++) {
    $sheet1.Cells.Item(2, $col).Formula = "=X_values*INPUT.FACTOR+INPUT.OFFSET+Control!`$C1-Control!`$D1"
    $sheet1.Cells.Item(3, $col).Formula = "=X_values^INPUT.FACTOR+INPUT.OFFSET+Control!`$C2-Control!`$D2"
}
